$wb = $excel.ActiveWorkbook

# Rename sheets
$wb.Worksheets.Item(1).Name = "GNG_TO-16502911267505567"
$wb.Worksheets.Item(2).Name = "NB_TO-16502911302325354"
$wb.Worksheets.Item(3).Name = "RS_TO-16502911302325354"
$wb.Worksheets.Item(4).Name = "TOL_TO-16502911302949843"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16502911303757212"

# Sheet 1 (GNG)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-1650291126718657.csv"
$ws1.Range("B3").Value = "GNG_stims-16502911267331543.csv"
$ws1.Range("B4").Value = "go_stims-16502911267331543.csv"
$ws1.Range("B5").Value = "GNG_stims-16502911267495275.csv"

# Sheet 2 (NB)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-16502911297820344.csv"
$ws2.Range("B3").Value = "TB-16502911300091178.csv"
$ws2.Range("B4").Value = "ZB-match_1-16502911267583895.csv"
$ws2.Range("B5").Value = "TB-16502911302119374.csv"
$ws2.Range("B6").Value = "OB-16502911274411602.csv"
$ws2.Range("B7").Value = "OB-16502911278682854.csv"
$ws2.Range("B8").Value = "ZB-match_7-16502911267998283.csv"
$ws2.Range("B9").Value = "OB-16502911273425515.csv"
$ws2.Range("B10").Value = "ZB-match_2-16502911270091844.csv"

# Sheet 4 (TOL)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16502911302478347.csv"
$ws4.Range("B3").Value = "ZM_stims-16502911302361298.csv"
$ws4.Range("B4").Value = "MM_stims-16502911302780163.csv"
$ws4.Range("B5").Value = "ZM_stims-16502911302488623.csv"
$ws4.Range("B6").Value = "MM_stims-16502911302940218.csv"
$ws4.Range("B7").Value = "ZM_stims-16502911302789814.csv"

# Sheet 5 (vSAT)
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-1650291130327612.csv"
$ws5.Range("B3").Value = "SAT_stims-16502911303119838.csv"
$ws5.Range("B4").Value = "vSAT_stims-16502911303599496.csv"
$ws5.Range("B5").Value = "SAT_stims-16502911302976038.csv"
